$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Season results" sheet: extend the totals formulas to also add in the
#    new "Week 8 results" sheet, and fix the stray 'Week 7 results'!C6 typo
#    in the C7 formula (it should reference C7, matching every other row).
# ---------------------------------------------------------------------------
$season = $wb.Worksheets.Item("Season results")
$season.Activate()

for ($row = 2; $row -le 7; $row++) {
    $season.Range("B$row").Formula = "='Week 4 results'!B$row+'Week 5 results'!B$row+'Week 6 results'!B$row+'Week 7 results'!B$row+'Week 8 results'!B$row"
    $season.Range("C$row").Formula = "='Week 4 results'!C$row+'Week 5 results'!C$row+'Week 6 results'!C$row+'Week 7 results'!C$row+'Week 8 results'!C$row"
    $season.Range("D$row").Formula = "='Week 4 results'!D$row+'Week 5 results'!D$row+'Week 6 results'!D$row+'Week 7 results'!D$row+'Week 8 results'!D$row"
}

$season.Range("B4").Select()

# ---------------------------------------------------------------------------
# 2) "Week 7 results" sheet: only the selected cell changes.
# ---------------------------------------------------------------------------
$week7 = $wb.Worksheets.Item("Week 7 results")
$week7.Activate()
$week7.Range("E2").Select()

# ---------------------------------------------------------------------------
# 3) "Week 8 results" sheet: fill in the week 8 matchups/results.
# ---------------------------------------------------------------------------
$week8 = $wb.Worksheets.Item("Week 8 results")
$week8.Activate()

# One of the models (AFNNET_12H2L_logsig, row 4) got a correct pick this week.
$week8.Range("B4").Value = 1

# Header row: matchup labels for each of the 12 games tested this week.
$week8.Range("F1").Value = "NYG(a) @ PHI(h) 2013 week 8"
$week8.Range("G1").Value = "SF(a) @ JAC(h) 2013 week 8"
$week8.Range("H1").Value = "DAL(a) @ DET(h) 2013 week 8"
$week8.Range("I1").Value = "CLE(a) @ KC(h) 2013 week 8"
$week8.Range("J1").Value = "MIA(a) @ NE(h) 2013 week 8"
$week8.Range("K1").Value = "BUF(a) @ NO(h) 2013 week 8"
$week8.Range("L1").Value = "PIT(a) @ OAK(h) 2013 week 8"
$week8.Range("M1").Value = "NYJ(a) @ CIN(h) 2013 week 8"
$week8.Range("N1").Value = "WAS(a) @ DEN(h) 2013 week 8"
$week8.Range("O1").Value = "ATL(a) @ ARI(h) 2013 week 8"
$week8.Range("P1").Value = "GB(a) @ MIN(h) 2013 week 8"
$week8.Range("Q1").Value = "SEA(a) @ STL(h) 2013 week 8"

# Row 4: the model's picks/probabilities for every game this week.
# (Set these before E1's rich-text value below, so new shared strings are
# created in the same order as the original edit.)
$week8.Range("E4").Value = "ON001(home win): 0.0004 / ON002(away win): 0.9995"
$week8.Range("F4").Value = "ON001(home win): 0.7312 / ON002(away win): 0.2687"
$week8.Range("G4").Value = "ON001(home win): 0.0001 / ON002(away win): 0.9998"
$week8.Range("H4").Value = "ON001(home win): 0.0215 / ON002(away win): 0.9784"
$week8.Range("I4").Value = "ON001(home win): 0.9997 / ON002(away win): 0.0002"
$week8.Range("J4").Value = "ON001(home win): 0.3261 / ON002(away win): 0.6738"
$week8.Range("K4").Value = "ON001(home win): 0.9983 / ON002(away win): 0.0016"
$week8.Range("L4").Value = "ON001(home win): 0.0125 / ON002(away win): 0.9874"
$week8.Range("M4").Value = "ON001(home win): 0.9615 / ON002(away win): 0.0384"
$week8.Range("N4").Value = "ON001(home win): 0.7056 / ON002(away win): 0.2943"
$week8.Range("O4").Value = "ON001(home win): 0.0016 / ON002(away win): 0.9983"
$week8.Range("P4").Value = "ON001(home win): 0.000004 / ON002(away win): 0.9999"
$week8.Range("Q4").Value = "ON001(home win): 0.0000009 / ON002(away win): 0.9999"

# E1 is the one matchup with the winning team (CAR) highlighted in red.
$week8.Range("E1").Value = "CAR(a) @ TB(h) 2013 week 8"
$e1 = $week8.Range("E1").Characters(1, 6)
$e1.Font.Bold = $true
$e1.Font.Color = 255
$rest = $week8.Range("E1").Characters(7, 20)
$rest.Font.Bold = $true

# Header row formatting: bold, centered horizontally and vertically,
# spanning one extra (blank) column same as the source edit.
$header = $week8.Range("A1:R1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# Column widths for the two trailing columns touched by the new layout.
$week8.Columns.Item(18).ColumnWidth = 11.166666666666666
$week8.Columns.Item(19).ColumnWidth = 9.307291666666666

# Match the printed page orientation used by the other weekly sheets.
$week8.PageSetup.Orientation = 1

$week8.Range("A9").Select()
$week8.Activate()
